$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "60.873.03"
$ws.Cells.Item(2, 5).Value = "  +1.03%  "
$ws.Cells.Item(3, 4).Value = "2.634.50"
$ws.Cells.Item(3, 5).Value = "  +1.85%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'529.10"
$ws.Cells.Item(5, 5).Value = "  +4.12%  "
$ws.Cells.Item(6, 4).Value = "'155.06"
$ws.Cells.Item(6, 5).Value = "  +1.25%  "
$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 5).Value = "  -0.06%  "
$ws.Cells.Item(8, 5).Value = "  -0.06%  "
$ws.Cells.Item(9, 5).Value = "  -0.16%  "
$ws.Cells.Item(10, 5).Value = "  +5.31%  "
$ws.Cells.Item(11, 5).Value = "  +1.50%  "
$ws.Cells.Item(12, 5).Value = "  +0.00%  "
$ws.Cells.Item(13, 4).Value = "3.095.70"
$ws.Cells.Item(13, 5).Value = "  +1.83%  "
$ws.Cells.Item(14, 4).Value = "60.875.85"
$ws.Cells.Item(14, 5).Value = "  +1.09%  "
$ws.Cells.Item(15, 4).Value = "'22.01"
$ws.Cells.Item(15, 5).Value = "  +2.55%  "
$ws.Cells.Item(16, 5).Value = "  +3.22%  "
$ws.Cells.Item(17, 4).Value = "2.639.08"
$ws.Cells.Item(17, 5).Value = "  +1.81%  "
$ws.Cells.Item(18, 5).Value = "  +0.48%  "
$ws.Cells.Item(19, 4).Value = "'353.08"
$ws.Cells.Item(19, 5).Value = "  +0.06%  "
$ws.Cells.Item(20, 4).Value = "'10.61"
$ws.Cells.Item(20, 5).Value = "  +1.25%  "
$ws.Cells.Item(21, 5).Value = "  +2.15%  "
$ws.Cells.Item(22, 5).Value = "  +0.34%  "
$ws.Cells.Item(23, 4).Value = "'61.59"
$ws.Cells.Item(23, 5).Value = "  +2.18%  "
$ws.Cells.Item(24, 5).Value = "  +2.58%  "
$ws.Cells.Item(25, 5).Value = "  +1.42%  "
$ws.Cells.Item(26, 4).Value = "'0.999"
$ws.Cells.Item(26, 5).Value = "  +0.12%  "
$ws.Cells.Item(27, 5).Value = "  +3.78%  "
$ws.Cells.Item(28, 5).Value = "  +1.14%  "
$ws.Cells.Item(29, 5).Value = "  -0.09%  "
$ws.Cells.Item(30, 4).Value = "'6.15"
$ws.Cells.Item(30, 5).Value = "  +7.81%  "
$ws.Cells.Item(31, 4).Value = "'19.47"
$ws.Cells.Item(31, 5).Value = "  +0.69%  "
$ws.Cells.Item(32, 5).Value = "  +4.13%  "
$ws.Cells.Item(33, 4).Value = "'150.33"
$ws.Cells.Item(33, 5).Value = "  -0.92%  "
$ws.Cells.Item(34, 5).Value = "  +4.75%  "
$ws.Cells.Item(35, 5).Value = "  +2.00%  "
$ws.Cells.Item(36, 5).Value = "  +10.68%  "
$ws.Cells.Item(37, 4).Value = "'0.890"
$ws.Cells.Item(37, 5).Value = "  +3.12%  "
$ws.Cells.Item(38, 5).Value = "  +1.63%  "
$ws.Cells.Item(39, 4).Value = "'3.82"
$ws.Cells.Item(39, 5).Value = "  +2.06%  "
$ws.Cells.Item(40, 4).Value = "'304.88"
$ws.Cells.Item(40, 5).Value = "  +3.61%  "
$ws.Cells.Item(41, 2).Value = "OKB"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(41, 4).Value = "'36.58"
$ws.Cells.Item(41, 5).Value = "  +1.46%  "
$ws.Cells.Item(42, 2).Value = "Mantle"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(42, 4).Value = "'0.640"
$ws.Cells.Item(42, 5).Value = "  +4.01%  "
$ws.Cells.Item(43, 2).Value = "Stellar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(43, 4).Value = "'0.102"
$ws.Cells.Item(43, 5).Value = "  +1.73%  "
$ws.Cells.Item(44, 2).Value = "Hedera"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(44, 4).Value = "'0.0562"
$ws.Cells.Item(44, 5).Value = "  +2.17%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).Value = "'0.998"
$ws.Cells.Item(45, 5).Value = "  +0.05%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "'19.74"
$ws.Cells.Item(46, 5).Value = "  +0.87%  "
$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(47, 4).Value = "'4.95"
$ws.Cells.Item(47, 5).Value = "  +3.71%  "
$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).Value = "'0.0238"
$ws.Cells.Item(48, 5).Value = "  +2.48%  "
$ws.Cells.Item(49, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(49, 4).Value = "'19.30"
$ws.Cells.Item(49, 5).Value = "  +8.72%  "
$ws.Cells.Item(50, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(50, 4).Value = "'10.34"
$ws.Cells.Item(50, 5).Value = "  +0.28%  "
$ws.Cells.Item(51, 2).Value = "Maker"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(51, 4).Value = "1.979.26"
$ws.Cells.Item(51, 5).Value = "  -0.34%  "
